$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.06874679424801
$ws.Cells.Item(2, 3).Value = 2.080277099878213
$ws.Cells.Item(2, 4).Value = -18.06874679424801
$ws.Cells.Item(2, 5).Value = -18.06874679424801
$ws.Cells.Item(2, 6).Value = -18.06874679424801
$ws.Cells.Item(2, 7).Value = -18.06874679424801
$ws.Cells.Item(2, 8).Value = -18.06874679424801
$ws.Cells.Item(2, 9).Value = -18.06874679424801
$ws.Cells.Item(2, 10).Value = -18.06874679424801
$ws.Cells.Item(2, 11).Value = -18.06874679424801

$ws.Cells.Item(3, 2).Value = -18.06874679424801
$ws.Cells.Item(3, 3).Value = -18.06874679424801
$ws.Cells.Item(3, 4).Value = -18.06874679424801
$ws.Cells.Item(3, 5).Value = -18.06874679424801
$ws.Cells.Item(3, 6).Value = -18.06874679424801
$ws.Cells.Item(3, 7).Value = -18.06874679424801
$ws.Cells.Item(3, 8).Value = -18.06874679424801
$ws.Cells.Item(3, 9).Value = 2.719162937035599
$ws.Cells.Item(3, 10).Value = -18.06874679424801
$ws.Cells.Item(3, 11).Value = -18.06874679424801

$ws.Cells.Item(4, 2).Value = -18.06874679424801
$ws.Cells.Item(4, 3).Value = 2.223468101130213
$ws.Cells.Item(4, 4).Value = 2.070516008274183
$ws.Cells.Item(4, 5).Value = -18.06874679424801
$ws.Cells.Item(4, 6).Value = 3.516038209459038
$ws.Cells.Item(4, 7).Value = -18.06874679424801
$ws.Cells.Item(4, 8).Value = 1.558140391420307
$ws.Cells.Item(4, 9).Value = -18.06874679424801
$ws.Cells.Item(4, 10).Value = 2.208559905504982
$ws.Cells.Item(4, 11).Value = -18.06874679424801

$ws.Cells.Item(5, 2).Value = -18.06874679424801
$ws.Cells.Item(5, 3).Value = 1.736337349981925
$ws.Cells.Item(5, 4).Value = -18.06874679424801
$ws.Cells.Item(5, 5).Value = -18.06874679424801
$ws.Cells.Item(5, 6).Value = -18.06874679424801
$ws.Cells.Item(5, 7).Value = 2.880416635130659
$ws.Cells.Item(5, 8).Value = -18.06874679424801
$ws.Cells.Item(5, 9).Value = -18.06874679424801
$ws.Cells.Item(5, 10).Value = -18.06874679424801
$ws.Cells.Item(5, 11).Value = -18.06874679424801

$ws.Cells.Item(6, 2).Value = -18.06874679424801
$ws.Cells.Item(6, 3).Value = -18.06874679424801
$ws.Cells.Item(6, 4).Value = -18.06874679424801
$ws.Cells.Item(6, 5).Value = -18.06874679424801
$ws.Cells.Item(6, 6).Value = -18.06874679424801
$ws.Cells.Item(6, 7).Value = -18.06874679424801
$ws.Cells.Item(6, 8).Value = -18.06874679424801
$ws.Cells.Item(6, 9).Value = -18.06874679424801
$ws.Cells.Item(6, 10).Value = -18.06874679424801
$ws.Cells.Item(6, 11).Value = -18.06874679424801

$ws.Cells.Item(7, 2).Value = 4.321923109899267
$ws.Cells.Item(7, 3).Value = -18.06874679424801
$ws.Cells.Item(7, 4).Value = -18.06874679424801
$ws.Cells.Item(7, 5).Value = -18.06874679424801
$ws.Cells.Item(7, 6).Value = -18.06874679424801
$ws.Cells.Item(7, 7).Value = -18.06874679424801
$ws.Cells.Item(7, 8).Value = -18.06874679424801
$ws.Cells.Item(7, 9).Value = -18.06874679424801
$ws.Cells.Item(7, 10).Value = -18.06874679424801
$ws.Cells.Item(7, 11).Value = -18.06874679424801

$ws.Cells.Item(8, 2).Value = -18.06874679424801
$ws.Cells.Item(8, 3).Value = -18.06874679424801
$ws.Cells.Item(8, 4).Value = -18.06874679424801
$ws.Cells.Item(8, 5).Value = 1.825320455485083
$ws.Cells.Item(8, 6).Value = -18.06874679424801
$ws.Cells.Item(8, 7).Value = -18.06874679424801
$ws.Cells.Item(8, 8).Value = -18.06874679424801
$ws.Cells.Item(8, 9).Value = -18.06874679424801
$ws.Cells.Item(8, 10).Value = -18.06874679424801
$ws.Cells.Item(8, 11).Value = -18.06874679424801

$ws.Cells.Item(9, 2).Value = -18.06874679424801
$ws.Cells.Item(9, 3).Value = -18.06874679424801
$ws.Cells.Item(9, 4).Value = -18.06874679424801
$ws.Cells.Item(9, 5).Value = -18.06874679424801
$ws.Cells.Item(9, 6).Value = -18.06874679424801
$ws.Cells.Item(9, 7).Value = -18.06874679424801
$ws.Cells.Item(9, 8).Value = -18.06874679424801
$ws.Cells.Item(9, 9).Value = -18.06874679424801
$ws.Cells.Item(9, 10).Value = -18.06874679424801
$ws.Cells.Item(9, 11).Value = -18.06874679424801

$ws.Cells.Item(10, 2).Value = -18.06874679424801
$ws.Cells.Item(10, 3).Value = -18.06874679424801
$ws.Cells.Item(10, 4).Value = -18.06874679424801
$ws.Cells.Item(10, 5).Value = -18.06874679424801
$ws.Cells.Item(10, 6).Value = -18.06874679424801
$ws.Cells.Item(10, 7).Value = -18.06874679424801
$ws.Cells.Item(10, 8).Value = -18.06874679424801
$ws.Cells.Item(10, 9).Value = 1.132766943255462
$ws.Cells.Item(10, 10).Value = -18.06874679424801
$ws.Cells.Item(10, 11).Value = 1.920044664069563

$ws.Cells.Item(11, 2).Value = -18.06874679424801
$ws.Cells.Item(11, 3).Value = -18.06874679424801
$ws.Cells.Item(11, 4).Value = -18.06874679424801
$ws.Cells.Item(11, 5).Value = 2.786663891382486
$ws.Cells.Item(11, 6).Value = -18.06874679424801
$ws.Cells.Item(11, 7).Value = 2.835720986168793
$ws.Cells.Item(11, 8).Value = -18.06874679424801
$ws.Cells.Item(11, 9).Value = -18.06874679424801
$ws.Cells.Item(11, 10).Value = -18.06874679424801
$ws.Cells.Item(11, 11).Value = 2.002240188768901

$ws.Cells.Item(12, 2).Value = -18.06874679424801
$ws.Cells.Item(12, 3).Value = -18.06874679424801
$ws.Cells.Item(12, 4).Value = -18.06874679424801
$ws.Cells.Item(12, 5).Value = -18.06874679424801
$ws.Cells.Item(12, 6).Value = -18.06874679424801
$ws.Cells.Item(12, 7).Value = -18.06874679424801
$ws.Cells.Item(12, 8).Value = -18.06874679424801
$ws.Cells.Item(12, 9).Value = -18.06874679424801
$ws.Cells.Item(12, 10).Value = -18.06874679424801
$ws.Cells.Item(12, 11).Value = -18.06874679424801

$ws.Cells.Item(13, 2).Value = -18.06874679424801
$ws.Cells.Item(13, 3).Value = -18.06874679424801
$ws.Cells.Item(13, 4).Value = -18.06874679424801
$ws.Cells.Item(13, 5).Value = 2.416871264955344
$ws.Cells.Item(13, 6).Value = -18.06874679424801
$ws.Cells.Item(13, 7).Value = -18.06874679424801
$ws.Cells.Item(13, 8).Value = -18.06874679424801
$ws.Cells.Item(13, 9).Value = -18.06874679424801
$ws.Cells.Item(13, 10).Value = 1.943862741924684
$ws.Cells.Item(13, 11).Value = 1.877035648439936

$ws.Cells.Item(14, 2).Value = -18.06874679424801
$ws.Cells.Item(14, 3).Value = -18.06874679424801
$ws.Cells.Item(14, 4).Value = 1.40626024651063
$ws.Cells.Item(14, 5).Value = -18.06874679424801
$ws.Cells.Item(14, 6).Value = -18.06874679424801
$ws.Cells.Item(14, 7).Value = -18.06874679424801
$ws.Cells.Item(14, 8).Value = -18.06874679424801
$ws.Cells.Item(14, 9).Value = -18.06874679424801
$ws.Cells.Item(14, 10).Value = -18.06874679424801
$ws.Cells.Item(14, 11).Value = 2.083895795318084

$ws.Cells.Item(15, 2).Value = -18.06874679424801
$ws.Cells.Item(15, 3).Value = -18.06874679424801
$ws.Cells.Item(15, 4).Value = 1.245103888882463
$ws.Cells.Item(15, 5).Value = -18.06874679424801
$ws.Cells.Item(15, 6).Value = -18.06874679424801
$ws.Cells.Item(15, 7).Value = -18.06874679424801
$ws.Cells.Item(15, 8).Value = -18.06874679424801
$ws.Cells.Item(15, 9).Value = -18.06874679424801
$ws.Cells.Item(15, 10).Value = -18.06874679424801
$ws.Cells.Item(15, 11).Value = -18.06874679424801

$ws.Cells.Item(16, 2).Value = -18.06874679424801
$ws.Cells.Item(16, 3).Value = -18.06874679424801
$ws.Cells.Item(16, 4).Value = -18.06874679424801
$ws.Cells.Item(16, 5).Value = -18.06874679424801
$ws.Cells.Item(16, 6).Value = -18.06874679424801
$ws.Cells.Item(16, 7).Value = -18.06874679424801
$ws.Cells.Item(16, 8).Value = -18.06874679424801
$ws.Cells.Item(16, 9).Value = -18.06874679424801
$ws.Cells.Item(16, 10).Value = 2.117441213188277
$ws.Cells.Item(16, 11).Value = -18.06874679424801

$ws.Cells.Item(17, 2).Value = -18.06874679424801
$ws.Cells.Item(17, 3).Value = 2.044345461321916
$ws.Cells.Item(17, 4).Value = 2.380854275146588
$ws.Cells.Item(17, 5).Value = -18.06874679424801
$ws.Cells.Item(17, 6).Value = -18.06874679424801
$ws.Cells.Item(17, 7).Value = -18.06874679424801
$ws.Cells.Item(17, 8).Value = 1.311724840058633
$ws.Cells.Item(17, 9).Value = 2.076748818158357
$ws.Cells.Item(17, 10).Value = 2.098721580410998
$ws.Cells.Item(17, 11).Value = -18.06874679424801

$ws.Cells.Item(18, 2).Value = -18.06874679424801
$ws.Cells.Item(18, 3).Value = -18.06874679424801
$ws.Cells.Item(18, 4).Value = -18.06874679424801
$ws.Cells.Item(18, 5).Value = -18.06874679424801
$ws.Cells.Item(18, 6).Value = -18.06874679424801
$ws.Cells.Item(18, 7).Value = -18.06874679424801
$ws.Cells.Item(18, 8).Value = 1.730530083556095
$ws.Cells.Item(18, 9).Value = 1.168969668085803
$ws.Cells.Item(18, 10).Value = 1.539993518389894
$ws.Cells.Item(18, 11).Value = -18.06874679424801

$ws.Cells.Item(19, 2).Value = -18.06874679424801
$ws.Cells.Item(19, 3).Value = -18.06874679424801
$ws.Cells.Item(19, 4).Value = 1.505478991151828
$ws.Cells.Item(19, 5).Value = -18.06874679424801
$ws.Cells.Item(19, 6).Value = -18.06874679424801
$ws.Cells.Item(19, 7).Value = -18.06874679424801
$ws.Cells.Item(19, 8).Value = 1.515020964379941
$ws.Cells.Item(19, 9).Value = 1.472200415153179
$ws.Cells.Item(19, 10).Value = -18.06874679424801
$ws.Cells.Item(19, 11).Value = -18.06874679424801

$ws.Cells.Item(20, 2).Value = -18.06874679424801
$ws.Cells.Item(20, 3).Value = 0.8277777843162328
$ws.Cells.Item(20, 4).Value = 1.449477967005355
$ws.Cells.Item(20, 5).Value = -18.06874679424801
$ws.Cells.Item(20, 6).Value = 3.09756893984621
$ws.Cells.Item(20, 7).Value = -18.06874679424801
$ws.Cells.Item(20, 8).Value = 1.999276996839708
$ws.Cells.Item(20, 9).Value = 0.986040328960617
$ws.Cells.Item(20, 10).Value = -18.06874679424801
$ws.Cells.Item(20, 11).Value = 2.103222816172557

$ws.Cells.Item(21, 2).Value = -18.06874679424801
$ws.Cells.Item(21, 3).Value = 0.9024402853905517
$ws.Cells.Item(21, 4).Value = -18.06874679424801
$ws.Cells.Item(21, 5).Value = 2.07575337632173
$ws.Cells.Item(21, 6).Value = -18.06874679424801
$ws.Cells.Item(21, 7).Value = 2.458727348228497
$ws.Cells.Item(21, 8).Value = 2.136947436344143
$ws.Cells.Item(21, 9).Value = -18.06874679424801
$ws.Cells.Item(21, 10).Value = -18.06874679424801
$ws.Cells.Item(21, 11).Value = -18.06874679424801

Write-Host "Done updating PSSM values."
